$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20: genomeweb.com link for the new Veracyte article
$link1 = "https://www.genomeweb.com/cancer/veracyte-eyes-dropping-ivdpharma-subsidiary-europe-growth-concentrates-among-clinical-ldts"
$ws.Range("A20").Value = $link1
$ws.Hyperlinks.Add($ws.Range("A20"), $link1) | Out-Null
$ws.Range("A20").Style = "Hyperlink"
$ws.Range("B20").Value = "Veracyte"
$ws.Range("C20").Value = "Veracyte Eyes Dropping IVD/Pharma Subsidiary in Europe As Growth Concentrates Among Clinical LDTs"

# Row 21: 360dx.com link for the same article
$link2 = "https://www.360dx.com/cancer/veracyte-eyes-dropping-ivdpharma-subsidiary-europe-growth-concentrates-among-clinical-ldts"
$ws.Range("A21").Value = $link2
$ws.Hyperlinks.Add($ws.Range("A21"), $link2) | Out-Null
$ws.Range("A21").Style = "Hyperlink"
$ws.Range("B21").Value = "Veracyte"
$ws.Range("C21").Value = "Veracyte Eyes Dropping IVD/Pharma Subsidiary in Europe As Growth Concentrates Among Clinical LDTs"

$wb.Save()
